# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.713.40"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.599.96"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "'211.41"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").Value = "'1.01"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").Value = "'0.0618"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "'0.247"
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").Value = "'19.54"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").Value = "'0.0841"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "1.825.25"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").Value = "1.593.44"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").Value = "26.686.76"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "0.0₃0756"
$ws.Range("E18").Value = "  +3.36%  "
$ws.Range("D19").Value = "'7.22"
$ws.Range("E19").Value = "  +3.69%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "'209.11"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'4.28"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").Value = "'2.30"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").Value = "'142.54"
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("D30").Value = "'0.0520"
$ws.Range("E30").Value = "  +2.52%  "
$ws.Range("D31").Value = "'1.15"
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("D32").Value = "'3.25"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").Value = "'2.97"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("D34").Value = "1.292.53"
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("E35").Value = "  -5.19%  "
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  +19.94%  "
$ws.Range("D40").Value = "'0.826"
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.19"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.785"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "'63.26"
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("D45").Value = "1.737.06"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'91.28"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "'7.38"
$ws.Range("E51").Value = "  -1.54%  "
